$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, matching the original inline-string cells
# (prices/percentages are stored as text like "240.30", not numbers, so we must
# stop Excel from auto-coercing the assignment to a Double and losing formatting).
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Refreshed prices / 1h-volume percentages (scraped data update), plus a handful of
# rows 47-51 that got re-ranked (coin name/link/price/volume swapped to new rows).
Set-TextValue "D2" '29.175.74'
Set-TextValue "D3" '1.834.23'
Set-TextValue "E3" '  +0.00%  '
Set-TextValue "D4" '0.9995'
Set-TextValue "E4" '  +0.05%  '
Set-TextValue "D5" '240.30'
Set-TextValue "E5" '  -2.01%  '
Set-TextValue "D6" '0.6852'
Set-TextValue "E6" '  -1.21%  '
Set-TextValue "D7" '0.9999'
Set-TextValue "E7" '  +0.02%  '
Set-TextValue "D8" '0.3016'
Set-TextValue "E8" '  -0.97%  '
Set-TextValue "D9" '0.07481'
Set-TextValue "E9" '  -2.67%  '
Set-TextValue "E10" '  -0.91%  '
Set-TextValue "D11" '0.07658'
Set-TextValue "E11" '  -1.93%  '
Set-TextValue "D12" '1.834.46'
Set-TextValue "E12" '  +0.07%  '
Set-TextValue "D13" '5.063'
Set-TextValue "E13" '  -0.77%  '
Set-TextValue "D14" '0.6826'
Set-TextValue "E14" '  +0.13%  '
Set-TextValue "D15" '87.12'
Set-TextValue "E15" '  -6.47%  '
Set-TextValue "D16" '6.198'
Set-TextValue "E16" '  -5.98%  '
Set-TextValue "D17" '29.169.07'
Set-TextValue "E17" '  +0.85%  '
Set-TextValue "D18" '0.000008182'
Set-TextValue "E18" '  -0.92%  '
Set-TextValue "D19" '2.082.58'
Set-TextValue "E19" '  +0.39%  '
Set-TextValue "D20" '12.55'
Set-TextValue "E20" '  -1.12%  '
Set-TextValue "D21" '226.33'
Set-TextValue "E21" '  -6.08%  '
Set-TextValue "D23" '7.424'
Set-TextValue "E23" '  -0.45%  '
Set-TextValue "D24" '1.001'
Set-TextValue "E24" '  +0.07%  '
Set-TextValue "D25" '0.1457'
Set-TextValue "E25" '  -3.09%  '
Set-TextValue "D26" '159.84'
Set-TextValue "E26" '  +0.94%  '
Set-TextValue "D27" '8.762'
Set-TextValue "E27" '  +0.07%  '
Set-TextValue "E28" '  -0.51%  '
Set-TextValue "D29" '1.502'
Set-TextValue "E29" '  -2.31%  '
Set-TextValue "D30" '4.260'
Set-TextValue "E30" '  +1.12%  '
Set-TextValue "D31" '4.141'
Set-TextValue "E31" '  -0.41%  '
Set-TextValue "D32" '1.207'
Set-TextValue "E32" '  +1.09%  '
Set-TextValue "D33" '0.05152'
Set-TextValue "E33" '  +0.86%  '
Set-TextValue "D34" '0.7673'
Set-TextValue "E34" '  -1.51%  '
Set-TextValue "D35" '1.842'
Set-TextValue "E35" '  -0.68%  '
Set-TextValue "E36" '  -1.06%  '
Set-TextValue "D37" '2.674'
Set-TextValue "E37" '  -0.84%  '
Set-TextValue "D38" '1.308.29'
Set-TextValue "E38" '  +1.44%  '
Set-TextValue "D39" '0.01832'
Set-TextValue "E39" '  -1.40%  '
Set-TextValue "D40" '2.715'
Set-TextValue "E40" '  +0.54%  '
Set-TextValue "D41" '0.9355'
Set-TextValue "E41" '  -1.93%  '
Set-TextValue "D42" '5.820'
Set-TextValue "E42" '  -5.20%  '
Set-TextValue "D43" '104.21'
Set-TextValue "E43" '  -2.44%  '
Set-TextValue "D44" '0.9997'
Set-TextValue "E44" '  +0.04%  '
Set-TextValue "D45" '65.22'
Set-TextValue "E45" '  +2.06%  '
Set-TextValue "D46" '9.596'
Set-TextValue "E46" '  -0.91%  '
Set-TextValue "B47" 'BabyDogeCoin'
Set-TextValue "C47" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D47" '0.00000000123'
Set-TextValue "E47" '  -0.16%  '
Set-TextValue "D48" '0.5203'
Set-TextValue "E48" '  +0.70%  '
Set-TextValue "B49" 'RocketPoolETH'
Set-TextValue "C49" 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue "D49" '1.981.00'
Set-TextValue "E49" '  +0.31%  '
Set-TextValue "B50" 'RenderToken'
Set-TextValue "C50" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D50" '1.770'
Set-TextValue "E50" '  +0.92%  '
Set-TextValue "B51" 'Cronos'
Set-TextValue "C51" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D51" '0.05918'
Set-TextValue "E51" '  +0.98%  '
